$p = $ppt.ActivePresentation

# --- 1) Update the cached "datetimeFigureOut" date text from 2021-12-17 to 2022-12-09
#        everywhere it appears (the slide master + every slide layout). ---
$oldDate = "2021-12-17"
$newDate = "2022-12-09"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShapes $layout.Shapes
}

# --- 2) Slide 1 text updates ---
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            $t = $shp.TextFrame.TextRange.Text
            if ($t -eq "Gephi Library") {
                $shp.TextFrame.TextRange.Text = "Gephi Library + JAVAFX"
            } elseif ($t -eq "Graph Driver") {
                $shp.TextFrame.TextRange.Text = "JDBC Driver For TurboGraph++"
            }
        }
    }
}
